$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.376.74'
$ws.Range('D3').Value = '1.937.40'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7714'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '248.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9987'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.10'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3215'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07096'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7883'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07996'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').Value = '1.933.32'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.387'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.07%  '
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.40%  '
$ws.Range('D17').Value = '30.375.71'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '258.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008034'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.815'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').Value = '2.193.36'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9984'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9994'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.824'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.623'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1359'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.307'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.371'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.530'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.446'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.176'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05218'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.293'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7532'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.769'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01978'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.815'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.478'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4530'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.989'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9994'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.580'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.837'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '986.52'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.56%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4182'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '
